$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.508.96'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +11.44%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.640.68'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +12.52%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '519.17'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +9.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '161.34'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +10.88%  '
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.613'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = 'USDC'
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.988'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.690.80'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +14.79%  '
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.107'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +11.90%  '
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = 'Toncoin'
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.17'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +13.47%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.351'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +7.78%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.56%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.081.96'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +11.95%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '60.966.65'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +10.51%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '22.62'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +13.33%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000143'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +10.60%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.652.68'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +13.01%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.86'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +6.85%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '364.01'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +15.88%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.70'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +11.87%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.25'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +10.69%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '61.09'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +8.77%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.430'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +9.23%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.171'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +13.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.704.46'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +10.83%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.984'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.53%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0882'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +18.97%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.66'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +8.63%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.995'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.55%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '19.93'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +10.28%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '158.20'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +8.58%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.60'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +8.43%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.66'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +11.43%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.23'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +12.50%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.04'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +12.40%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.887'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +9.68%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.52'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +14.77%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '311.08'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +24.80%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.81'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +12.97%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '35.71'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +6.30%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.824'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +33.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.647'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +12.51%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0582'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +13.29%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.86%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '20.18'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +21.32%  '
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.06'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +15.75%  '
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'FirstDigitalUSD'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.982'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.55%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0241'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +9.27%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.044.84'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +14.22%  '

Write-Host "Applied cryptos update"
